# Auto-generated edit script: refreshes "ランサーズ" sheet data with the
# 2025-10-03 12:33:52 JST scrape results (19 data rows, header in row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the hyperlink objects (and their relationships) created by the
# previous scrape before rebuilding the data grid.
$ws.Hyperlinks.Delete()

# Wipe the previous data body (row 1 header stays untouched) so no stale
# values linger in cells that end up blank in the refreshed data.
$ws.Range("A2:H20").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(2, 2).Value = '【業務効率化】chatgpt×Googleスプレッドシートを使って教育カリキュラムの作成依頼'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5405813'
$ws.Cells.Item(2, 7).Value = 398
$ws.Cells.Item(2, 8).Value = '🔥GPT,ChatGPT ◆効率化'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5405813') | Out-Null

# Row 3
$ws.Cells.Item(3, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(3, 2).Value = '【業務自動化×補助金対応】生成AI活用/日本人モデル画像生成歓迎'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5405834'
$ws.Cells.Item(3, 7).Value = 395
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◆自動化'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5405834') | Out-Null

# Row 4
$ws.Cells.Item(4, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(4, 2).Value = '詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5377709'
$ws.Cells.Item(4, 7).Value = 245
$ws.Cells.Item(4, 8).Value = '🔥Next.js ◆開発,Node.js ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5377709') | Out-Null

# Row 5
$ws.Cells.Item(5, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(5, 2).Value = '<Next.js、バックエンド開発> ガントチャートアプリの改修製造'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5379158'
$ws.Cells.Item(5, 7).Value = 225
$ws.Cells.Item(5, 8).Value = '🔥Next.js ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5379158') | Out-Null

# Row 6
$ws.Cells.Item(6, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(6, 2).Value = 'Reactの細かい修正の対応'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5405740'
$ws.Cells.Item(6, 7).Value = 120
$ws.Cells.Item(6, 8).Value = '🔥React'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5405740') | Out-Null

# Row 7
$ws.Cells.Item(7, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(7, 2).Value = '【急募】愛知県でのBtoB受発注システム開発者募集'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5405971'
$ws.Cells.Item(7, 7).Value = 118
$ws.Cells.Item(7, 8).Value = '◆開発,システム開発'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5405971') | Out-Null

# Row 8
$ws.Cells.Item(8, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(8, 2).Value = '【RPA/Power Automate】税務システム自動化プロジェクトの依頼'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5403634'
$ws.Cells.Item(8, 7).Value = 103
$ws.Cells.Item(8, 8).Value = '◆自動化'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5403634') | Out-Null

# Row 9
$ws.Cells.Item(9, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(9, 2).Value = '【急募】ガチャ型ECサイト開発、その他案件にごお協力いただけるフリーランス募集!'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5406171'
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = '◆開発 ◇サイト'
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5406171') | Out-Null

# Row 10
$ws.Cells.Item(10, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(10, 2).Value = 'Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5379176'
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = '◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5379176') | Out-Null

# Row 11
$ws.Cells.Item(11, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(11, 2).Value = '【急募】WordPress予約カレンダープラグイン開発の依頼'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5406144'
$ws.Cells.Item(11, 7).Value = 88
$ws.Cells.Item(11, 8).Value = '◆開発 ○WordPress'
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5406144') | Out-Null

# Row 12
$ws.Cells.Item(12, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(12, 2).Value = '【WEB】Nuxt3でのWEBページ表示速度改善、他継続して弊社システムの開発に携われる方'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5406001'
$ws.Cells.Item(12, 7).Value = 83
$ws.Cells.Item(12, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5406001') | Out-Null

# Row 13
$ws.Cells.Item(13, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(13, 2).Value = '【急募】Excelマクロでデータからグラフを自動作成するツール'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5405961'
$ws.Cells.Item(13, 7).Value = 68
$ws.Cells.Item(13, 8).Value = '◆ツール'
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5405961') | Out-Null

# Row 14
$ws.Cells.Item(14, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(14, 2).Value = '【急募】国内WiFiレンタルサービスのショッピングカート移行(ECサイト構築)'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5406225'
$ws.Cells.Item(14, 7).Value = 45
$ws.Cells.Item(14, 8).Value = '◇サイト'
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5406225') | Out-Null

# Row 15
$ws.Cells.Item(15, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(15, 2).Value = '【急募】WebRelease2制作サイトの改修依頼'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5406212'
$ws.Cells.Item(15, 7).Value = 38
$ws.Cells.Item(15, 8).Value = '◇サイト'
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5406212') | Out-Null

# Row 16
$ws.Cells.Item(16, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(16, 2).Value = '【急募】Appsheetでの貸出・返却システム構築方法を教えてください'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5406372'
$ws.Cells.Item(16, 7).Value = 25
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), 'https://www.lancers.jp/work/detail/5406372') | Out-Null

# Row 17
$ws.Cells.Item(17, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(17, 2).Value = '【急募】国内300店舗規模のスーパーマーケット向けActive Directory構築'
$ws.Cells.Item(17, 3).Value = 'システム開発'
$ws.Cells.Item(17, 4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(17, 5).Value = '期限情報なし'
$ws.Cells.Item(17, 6).Value = 'https://www.lancers.jp/work/detail/5406008'
$ws.Cells.Item(17, 7).Value = 25
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), 'https://www.lancers.jp/work/detail/5406008') | Out-Null

# Row 18
$ws.Cells.Item(18, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(18, 2).Value = '【急募】Scala技術者募集!Googleセーフブラウジング対応'
$ws.Cells.Item(18, 3).Value = 'システム開発'
$ws.Cells.Item(18, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(18, 5).Value = '期限情報なし'
$ws.Cells.Item(18, 6).Value = 'https://www.lancers.jp/work/detail/5406304'
$ws.Cells.Item(18, 7).Value = 18
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), 'https://www.lancers.jp/work/detail/5406304') | Out-Null

# Row 19
$ws.Cells.Item(19, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(19, 2).Value = 'サービスLPと受付フォーム(クレジット決済機能)、入力内容を固定フォームで表示の制作依頼、'
$ws.Cells.Item(19, 3).Value = 'システム開発'
$ws.Cells.Item(19, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(19, 5).Value = '期限情報なし'
$ws.Cells.Item(19, 6).Value = 'https://www.lancers.jp/work/detail/5406154'
$ws.Cells.Item(19, 7).Value = 18
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), 'https://www.lancers.jp/work/detail/5406154') | Out-Null

# Row 20
$ws.Cells.Item(20, 1).Value = '2025-10-03 12:33:52'
$ws.Cells.Item(20, 2).Value = '【急募】全国物件情報抽出プログラム作成依頼'
$ws.Cells.Item(20, 3).Value = 'システム開発'
$ws.Cells.Item(20, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(20, 5).Value = '期限情報なし'
$ws.Cells.Item(20, 6).Value = 'https://www.lancers.jp/work/detail/5405763'
$ws.Cells.Item(20, 7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), 'https://www.lancers.jp/work/detail/5405763') | Out-Null

# Widen the skill-summary column (H) now that entries are longer.
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668

Write-Host "Refreshed" $wb.Worksheets.Item(1).Name "with" (20-1) "data rows"

